$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The naive-component forecaster bug fix produced one additional data point,
# so the whole error-stats table (columns B:G, rows 2-11) shifts down by one
# row - each row now shows the stats for the *previous* forecast horizon -
# while column A (the fixed Q0..Q9 / horizon-index labels) stays put. The
# oldest row (row 11) falls off the bottom of the table and the freshly
# computed values are written into the now-vacant row 2.

# Shift existing rows down, starting from the bottom so we don't overwrite
# values before they've been read.
for ($r = 10; $r -ge 2; $r--) {
    $src = $ws.Range("B" + $r + ":G" + $r)
    $dst = $ws.Range("B" + ($r + 1) + ":G" + ($r + 1))
    $dst.Value2 = $src.Value2
}

# Write the newly computed error statistics into the freed-up top row.
$ws.Range("B2").Value2 = -0.01725120502155203
$ws.Range("C2").Value2 = 1.384110966838059
$ws.Range("D2").Value2 = 8.376137918165924
$ws.Range("E2").Value2 = 2.894155821334768
$ws.Range("F2").Value2 = 2.959148581664683
$ws.Range("G2").Value2 = 23
